# Refresh the cryptocurrency Price (column D) and Volume(1h) change
# (column E) figures with the latest values from the data source.
# Price values are kept as plain text (matching the existing sheet
# layout) so values such as "0.9994" or "134.50" are not reinterpreted
# as numbers and do not lose significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.955.79"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.873.56"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.27"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4955"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.49"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2904"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06569"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "1.872.47"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.67"
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07166"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6575"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.91"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.844"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "29.941.21"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007832"
$ws.Range("E18").Value = "  +3.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "2.113.87"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.568"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.026"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "149.76"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.50"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.65"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.895"
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.371"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.143"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08679"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.923"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04999"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6975"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.089"
$ws.Range("E36").Value = "  -4.89%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.688"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.160"
$ws.Range("E39").Value = "  -5.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01695"
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9255"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9978"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4166"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.42"
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.372"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05640"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.32"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.54"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.328"
$ws.Range("E51").Value = "  -0.75%  "
